# Refresh the cryptos list: updated prices / 1h volume deltas, plus the
# Cardano/Dogecoin rows swapping rank order (rows 8 and 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric need to be forced to text
# so Excel does not auto-convert them (matching the source data which
# stores these as inline strings, e.g. "4.817", "0.6086"). We set the
# cell to Text format just long enough to assign the literal value, then
# restore the "Normal" style so formatting matches the rest of the sheet.
function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '29.227.19'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.829.43'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.13%  '
Set-TextValue 'D5' '238.17'
$ws.Range('E5').Value = '  -0.93%  '
Set-TextValue 'D6' '0.6086'
$ws.Range('E6').Value = '  -3.39%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D8' '0.2828'
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D9' '0.07075'
$ws.Range('E9').Value = '  -5.02%  '
Set-TextValue 'D10' '23.98'
$ws.Range('E10').Value = '  -3.75%  '
Set-TextValue 'D11' '0.07643'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '1.829.12'
$ws.Range('E12').Value = '  -0.81%  '
Set-TextValue 'D13' '4.817'
Set-TextValue 'D14' '0.6378'
$ws.Range('E14').Value = '  -5.95%  '
Set-TextValue 'D15' '0.000009967'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('E16').Value = '  -0.91%  '
Set-TextValue 'D17' '79.81'
$ws.Range('E17').Value = '  -2.80%  '
Set-TextValue 'D18' '5.992'
$ws.Range('E18').Value = '  -4.43%  '
$ws.Range('D19').Value = '29.211.42'
$ws.Range('E19').Value = '  -0.44%  '
Set-TextValue 'D20' '230.65'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E21').Value = '  -4.06%  '
$ws.Range('E22').Value = '  +0.12%  '
Set-TextValue 'D23' '7.022'
$ws.Range('E23').Value = '  -5.26%  '
Set-TextValue 'D24' '1.001'
$ws.Range('E24').Value = '  +0.11%  '
Set-TextValue 'D25' '155.51'
$ws.Range('E25').Value = '  -1.68%  '
Set-TextValue 'D26' '8.075'
$ws.Range('E26').Value = '  -4.86%  '
$ws.Range('E27').Value = '  -4.39%  '
$ws.Range('E28').Value = '  -3.80%  '
Set-TextValue 'D29' '0.06703'
$ws.Range('E29').Value = '  +2.74%  '
Set-TextValue 'D30' '1.464'
$ws.Range('E30').Value = '  +1.25%  '
Set-TextValue 'D31' '1.462'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('E32').Value = '  -5.24%  '
Set-TextValue 'D33' '3.817'
$ws.Range('E33').Value = '  -6.26%  '
Set-TextValue 'D34' '1.130'
$ws.Range('E34').Value = '  -0.85%  '
Set-TextValue 'D35' '1.729'
$ws.Range('E35').Value = '  -5.92%  '
Set-TextValue 'D36' '0.6576'
$ws.Range('E36').Value = '  -5.63%  '
Set-TextValue 'D37' '2.549'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('D38').Value = '1.238.16'
$ws.Range('E38').Value = '  -0.74%  '
Set-TextValue 'D39' '2.754'
$ws.Range('E39').Value = '  -2.24%  '
Set-TextValue 'D40' '0.01768'
$ws.Range('E40').Value = '  -4.75%  '
Set-TextValue 'D41' '6.601'
$ws.Range('E41').Value = '  -2.59%  '
Set-TextValue 'D42' '0.9281'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '1.983.74'
$ws.Range('E44').Value = '  -0.51%  '
Set-TextValue 'D45' '100.58'
$ws.Range('E45').Value = '  -0.19%  '
Set-TextValue 'D46' '63.71'
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('E48').Value = '  -4.84%  '
Set-TextValue 'D49' '8.525'
$ws.Range('E49').Value = '  -5.21%  '
$ws.Range('E50').Value = '  -5.11%  '
Set-TextValue 'D51' '0.05578'
$ws.Range('E51').Value = '  -1.58%  '
